# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime columns for the four
# localized-but-not-yet-handed-off files (rows 4-7) on both the "zh-cn"
# and "de-de" status sheets, reflecting a fresh handoff-report generation.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E4").Value = "ht"
$ws_zhcn.Range("E5").Value = "ht"
$ws_zhcn.Range("E6").Value = "ht"
$ws_zhcn.Range("E7").Value = "ht"
$ws_zhcn.Range("H4").Value = "2016-08-24 20:32:27"
$ws_zhcn.Range("H5").Value = "2016-08-24 20:32:27"
$ws_zhcn.Range("H6").Value = "2016-08-24 20:32:27"
$ws_zhcn.Range("H7").Value = "2016-08-24 20:32:27"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E4").Value = "ht"
$ws_dede.Range("E5").Value = "ht"
$ws_dede.Range("E6").Value = "ht"
$ws_dede.Range("E7").Value = "ht"
$ws_dede.Range("H4").Value = "2016-08-24 20:32:33"
$ws_dede.Range("H5").Value = "2016-08-24 20:32:33"
$ws_dede.Range("H6").Value = "2016-08-24 20:32:33"
$ws_dede.Range("H7").Value = "2016-08-24 20:32:33"
